$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for cryptos list refresh.
# D-column values are numeric-looking text; prefix with an apostrophe so Excel
# stores them as text (matching the inlineStr source data) and reset the style
# afterwards so no stray quote-prefix / text-format style is left on the cell.

$ws.Range("D2").Value = "'29.189.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "'1.857.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'242.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'0.7000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.07812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "'0.3114"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "'24.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").Value = "'0.07803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("D12").Value = "'1.851.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "'5.132"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "'92.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "'0.6918"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "'6.588"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "'0.000008515"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'29.203.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'248.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "'2.109.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'0.9993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'7.573"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'160.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").Value = "'8.923"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "'18.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "'1.570"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("D30").Value = "'4.277"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "'4.242"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'1.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'0.05247"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'0.7585"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "'1.872"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").Value = "'1.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").Value = "'1.231.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "'2.736"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "'0.9015"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'110.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").Value = "'5.843"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.60%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -7.95%  "
$ws.Range("D46").Value = "'2.008.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Value = "'0.5178"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "'9.510"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'1.765"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("E51").Value = "  -2.21%  "
